# Update "resumen_pedido" so all sections appear: for every row where the
# "Diferencia Stock" (column M) was non-zero, fold that adjustment into the
# "Unidades Pedido" (column L) and zero out the difference.
#
#   new L = old L - old M
#   new M = 0
#
# Afterwards refresh the summary metrics block:
#   C100 (Total_Unidades)      -= sum of the old M values that were zeroed
#   C111 (Total_Ajuste_Stock)   = 0

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> new value of "Unidades Pedido" (column L). "Diferencia Stock"
# (column M) is always reset to 0 for these rows.
$updates = [ordered]@{
    4  = 18
    6  = 6
    14 = 8
    15 = 1
    17 = 4
    23 = 4
    26 = 20
    27 = 4
    30 = 2
    31 = 12
    34 = 1
    35 = 20
    36 = 1
    37 = 1
    39 = 5
    45 = 8
    47 = 5
    53 = 2
    59 = 6
    62 = 6
    63 = 18
    67 = 12
    68 = 17
    70 = 4
    71 = 1
    80 = 4
    81 = 6
    90 = 6
    92 = 4
    96 = 1
}

$totalAdjustRemoved = 0

foreach ($row in $updates.Keys) {
    $mCell = $ws.Range("M$row")
    $oldM = $mCell.Value2
    $totalAdjustRemoved += $oldM

    $ws.Range("L$row").Value = $updates[$row]
    $mCell.Value = 0
}

# Refresh the "METRICAS DE RESUMEN" summary block.
$totalUnidadesCell = $ws.Range("C100")
$totalUnidadesCell.Value = $totalUnidadesCell.Value2 - $totalAdjustRemoved

$ws.Range("C111").Value = 0
